$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.768.28'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +2.53%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.375.08'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.32%  '

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.31%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '593.90'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +7.06%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '186.94'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.73%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.602'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +3.74%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.13%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.184'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +3.89%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.591'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +1.96%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '47.51'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +3.12%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000279'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +4.80%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.922.49'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.30%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '641.52'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +10.08%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.63'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.63%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '67.779.55'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.39%  '

$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.41%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.375.21'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.89%  '

$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.58%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.12'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.90%  '

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.97%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '17.94'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.66%  '

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.99%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '99.95'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.12%  '

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +2.39%  '

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +6.67%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.75'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +4.18%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '32.96'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +7.69%  '

$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +3.65%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.92'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +4.80%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '612.53'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +6.28%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.80'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.37%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.053.23'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +7.99%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.12'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +2.27%  '

$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +3.03%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.998'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.06%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '56.36'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.16%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.79'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +6.18%  '

$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.131'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +3.57%  '

$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'InjectiveProtocol'
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '33.87'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.18%  '

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +3.24%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0₃0700'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.03%  '

$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.12%  '

$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +2.57%  '

$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.85%  '

$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.52%  '

$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.56%  '

$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +12.16%  '

$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.09%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '128.05'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.55%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.73'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +5.41%  '
